{"js": "// Locate the \"<m>manure\" text (unique in the document): the <m> tag run\n// immediately followed by the \"manure\" word (itself split across two runs,\n// \"m\" + \"anure\").\nconst body = context.document.body;\n\n// --- Step 1: split the \"<m>\" tag run into \"<m\" + \">\" --------------------\n// Scope the search to the unique \"<m>manure\" phrase first so we grab the\n// exact \"<m>\" run that precedes \"manure\" (there are many other \"<m>\" tags\n// in this document).\nconst scopeRange = body.search(\"<m>manure\", { matchCase: true }).getFirst();\nawait context.sync();\n\nconst tagRun = scopeRange.search(\"<m>\", { matchCase: true }).getFirst();\nawait context.sync();\ntagRun.load(\"text\");\nawait context.sync();\n\n// split(\"m\", ...) on \"<m>\" yields two ranges: \"<m\" and \">\" (delimiter \"m\"\n// stays attached to the first piece because trimDelimiters = false).\nconst tagParts = tagRun.split([\"m\"], false, false);\ntagParts.load(\"text\");\nawait context.sync();\n\nconst openPart = tagParts.items[0];   // \"<m\"  (keeps original Courier-New/blue/not-bold formatting)\nconst closePart = tagParts.items[1];  // \">\"   (same Courier-New/blue styling, but no explicit bold-off)\nclosePart.font.bold = false;\nawait context.sync();\n\n// --- Step 2: replace \"manure\" with \"earth mixed with dung\" --------------\nconst wordRange = body.search(\"manure\", { matchCase: true }).getFirst();\nawait context.sync();\nwordRange.insertText(\"earth mixed with dung\", Word.InsertLocation.replace);\nawait context.sync();\n\nconst newWordRange = body.search(\"earth mixed with dung\", { matchCase: true }).getFirst();\nawait context.sync();\nnewWordRange.font.bold = false;\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# \"<m>manure\" (the <m> tag immediately followed by the word \"manure\", itself\n# split across two runs \"m\"+\"anure\") is unique in the document, so Find\n# unambiguously locates the exact spot to edit.\n$rng = $d.Content\n$rng.Find.ClearFormatting()\n$rng.Find.MatchCase = $true\n$rng.Find.MatchWholeWord = $false\n$rng.Find.MatchWildcards = $false\n$found = $rng.Find.Execute(\"<m>manure\")\n\nif ($found) {\n    $start = $rng.Start\n    $end = $rng.End\n\n    # --- Split the \"<m>\" tag run into \"<m\" + \">\" -------------------------\n    # Characters: [start, start+2) = \"<m\", [start+2, start+3) = \">\".\n    # Flipping Bold off on just the \">\" character forces Word to break it\n    # into its own run (distinct formatting from the \"<m\" run, which keeps\n    # its explicit bold-off).\n    $gtRng = $d.Range($start + 2, $start + 3)\n    $gtRng.Font.Bold = 0\n\n    # --- Replace \"manure\" with \"earth mixed with dung\" -------------------\n    $wordRng = $d.Range($start + 3, $end)\n    $wordRng.Text = \"earth mixed with dung\"\n    $wordRng.Font.Bold = 0\n}\n"}
